$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A (former A shifts to B)
$ws.Columns.Item(1).Insert()

# New column A width (closest achievable to the authored 19.33203125 chars
# given this runtime's pixel-quantized ColumnWidth property)
$ws.Columns.Item(1).ColumnWidth = 18.5

# Header row: new "Mã nhóm hoạt động" header in A1, bold, sharing B1's
# fill/border/alignment (only the font differs from B1's header style)
$ws.Range("A1").Value = "Mã nhóm hoạt động"
$ws.Range("A1").Font.Bold = $true
$ws.Range("A1").Interior.Color = $ws.Range("B1").Interior.Color
$ws.Range("A1").Borders.Color = $ws.Range("B1").Borders.Color
$ws.Range("A1").Borders.LineStyle = $ws.Range("B1").Borders.LineStyle

# Data rows
$ws.Range("A2").Value = "Nhom001"
$ws.Range("B2").Value = "Thể dục thể thao"

$ws.Range("A3").Value = "Nhom002"
$ws.Range("B3").Value = "Đường lối định hướng"

# Header row height to match the authored value
$ws.Rows.Item(1).RowHeight = 15.6

# Restore selection to match the authored state
[void]$ws.Range("B12").Select()
